# "Add files via upload" — slide 5 ("IMPLEMENTATION AND CODE") has a
# reference table whose last row ("PPT Link") had an empty second
# (link) cell. Fill it in with the link to this very presentation on
# GitHub, then let the table's shape frame grow to accommodate the
# extra wrapped lines the long URL now needs.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shp = $s.Shapes.Item("Table 6")
$tbl = $shp.Table

# Row 4 = "PPT Link", Column 2 = the link cell (currently empty).
$cell = $tbl.Cell(4, 2)
$cell.Shape.TextFrame.TextRange.Text = "https://github.com/dharrini06/WeatherPredictionAI/blob/5b303c34095f5efe1d00307f42949aadb7eab5b6/Weather%20Prediction%20AI%20mini%20project.pptx"

# The newly-wrapped URL needs a taller last row, which grows the
# overall table frame (cx stays the same, cy grows) just like it does
# in PowerPoint when a cell's text no longer fits on one line.
$row4 = $tbl.Rows.Item(4)
$row4.Height = 208.0
